# Update feed logs and data lake files
# Append two new log rows (74, 75) to Sheet1, matching the existing data pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(74, 1).Value = 73
$ws.Cells.Item(74, 2).Value = 1
$ws.Cells.Item(74, 3).Value = "2024-06-16 11:11:00"
$ws.Cells.Item(74, 4).Value = 200
$ws.Cells.Item(74, 5).Value = 7

$ws.Cells.Item(75, 1).Value = 74
$ws.Cells.Item(75, 2).Value = 2
$ws.Cells.Item(75, 3).Value = "2024-06-16 11:11:00"
$ws.Cells.Item(75, 4).Value = 200
$ws.Cells.Item(75, 5).Value = 0
